# Update Portugal MSME country indicator figures (text values stored as strings).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D11" = "107.56"
    "D12" = "77.35"
    "B33" = "74.83"
    "C33" = "3.87"
    "B34" = "41.39"
    "C34" = "37.64"
    "D34" = "79.03"
    "B36" = "94.98"
    "C36" = "4.92"
    "B40" = "23.08"
    "C40" = "44.65"
    "D40" = "67.73"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # keep these as text values (matching the original shared-string cells)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Host "Updated cells: $($updates.Keys -join ', ')"
